$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.932
$ws.Range("E3").Value = 12.932
$ws.Range("E5").Value = 13.165
$ws.Range("C9").Value = -11.899
$ws.Range("E11").Value = 13.012
$ws.Range("E12").Value = 13
$ws.Range("C13").Value = -12.201
$ws.Range("C16").Value = -12.439
$ws.Range("C18").Value = -12.362
$ws.Range("C20").Value = -12.732
$ws.Range("E21").Value = 13.166
